# All output fields used in the combo_order_eval table must have values.
# Row 7 and row 8 were missing a value in column G (the "accumulate"
# input column), and row 8 was missing its FedEx-zone output values
# (columns I:L) that every other data row already has.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the previously-blank "accumulate" input cells with 0, matching
# the pattern used by the other data rows.
$ws.Range("G7").Value = 0
$ws.Range("G8").Value = 0

# Row 8 was missing handling/zone output values entirely; copy the
# formatting from the equivalent cells of row 6 (same "N/A" pattern used
# when a row has no FedEx Styrofoam Box zone pricing) and fill them in.
$ws.Range("I6:L6").Copy()
$ws.Range("I8:L8").PasteSpecial(-4122)

$ws.Range("I8").Value = "N/A"
$ws.Range("J8").Value = "N/A"
$ws.Range("K8").Value = "N/A"
$ws.Range("L8").Value = "N/A"

# Leave the selection where the author left it when finishing the edit.
$ws.Range("G12").Select()
